# The deck ships two embedded theme parts:
#   ppt/theme/theme1.xml  -> used by the Notes Master ("Office Theme" colours)
#   ppt/theme/theme2.xml  -> used by the Slide Master / all slides ("Integral" colours)
#
# The authored change swaps the two themes' content so the slides (and the
# presentation's default theme) end up using the plain "Office Theme" colour
# palette instead of the green/yellow/teal "Integral" palette.
#
# The scheme/font/format (fontScheme + fmtScheme) blocks of the two themes are
# already byte-identical, so the only thing that actually needs to change is
# the 12-slot theme colour scheme used by the slide master's theme part
# (ppt/theme/theme2.xml): it must move from the "Integral" palette to the
# "Office" palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $tcs.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# msoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink.
Set-ThemeColor 1  "000000"   # dk1
Set-ThemeColor 2  "FFFFFF"   # lt1
Set-ThemeColor 3  "44546A"   # dk2
Set-ThemeColor 4  "E7E6E6"   # lt2
Set-ThemeColor 5  "5B9BD5"   # accent1
Set-ThemeColor 6  "ED7D31"   # accent2
Set-ThemeColor 7  "A5A5A5"   # accent3
Set-ThemeColor 8  "FFC000"   # accent4
Set-ThemeColor 9  "4472C4"   # accent5
Set-ThemeColor 10 "70AD47"   # accent6
Set-ThemeColor 11 "0563C1"   # hlink
Set-ThemeColor 12 "954F72"   # folHlink
